# Apply crypto price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.590.14'
$ws.Range("E2").Value = '  +5.06%  '

$ws.Range("D3").Value = '3.524.93'
$ws.Range("E3").Value = '  +14.56%  '

$ws.Range("D4").Value = '''0.997'
$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").Value = '''590.11'
$ws.Range("E5").Value = '  +2.73%  '

$ws.Range("D6").Value = '''185.15'
$ws.Range("E6").Value = '  +8.86%  '

$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").Value = '3.534.12'
$ws.Range("E7").Value = '  +14.95%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '''0.996'
$ws.Range("E8").Value = '  -0.45%  '

$ws.Range("D9").Value = '''0.533'
$ws.Range("E9").Value = '  +4.69%  '

$ws.Range("D10").Value = '''6.59'
$ws.Range("E10").Value = '  +4.29%  '

$ws.Range("D11").Value = '''0.158'
$ws.Range("E11").Value = '  +5.89%  '

$ws.Range("D12").Value = '''0.490'
$ws.Range("E12").Value = '  +4.19%  '

$ws.Range("D13").Value = '''38.69'
$ws.Range("E13").Value = '  +8.20%  '

$ws.Range("D14").Value = '''0.0000250'
$ws.Range("E14").Value = '  +4.56%  '

$ws.Range("D15").Value = '4.042.62'
$ws.Range("E15").Value = '  +12.52%  '

$ws.Range("D16").Value = '69.752.41'
$ws.Range("E16").Value = '  +5.35%  '

$ws.Range("E17").Value = '  +1.28%  '

$ws.Range("D18").Value = '3.451.80'
$ws.Range("E18").Value = '  +12.08%  '

$ws.Range("D19").Value = '''7.43'
$ws.Range("E19").Value = '  +6.97%  '

$ws.Range("D20").Value = '''17.06'
$ws.Range("E20").Value = '  +2.84%  '

$ws.Range("D21").Value = '''503.57'
$ws.Range("E21").Value = '  +3.42%  '

$ws.Range("D22").Value = '''8.84'
$ws.Range("E22").Value = '  +15.08%  '

$ws.Range("D23").Value = '''0.732'
$ws.Range("E23").Value = '  +6.75%  '

$ws.Range("D24").Value = '''86.60'
$ws.Range("E24").Value = '  +4.92%  '

$ws.Range("D25").Value = '''13.32'
$ws.Range("E25").Value = '  +5.33%  '

$ws.Range("D26").Value = '''2.38'
$ws.Range("E26").Value = '  +8.06%  '

$ws.Range("E27").Value = '  +4.93%  '

$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").Value = '''2.52'
$ws.Range("E29").Value = '  +12.24%  '

$ws.Range("D30").Value = '''8.13'
$ws.Range("E30").Value = '  +3.32%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '''30.97'
$ws.Range("E31").Value = '  +11.71%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''2.72'
$ws.Range("E32").Value = '  +4.83%  '

$ws.Range("B33").Value = 'PEPE'
$ws.Range("C33").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D33").Value = '''0.0000107'
$ws.Range("E33").Value = '  +19.29%  '

$ws.Range("D34").Value = '''0.117'
$ws.Range("E34").Value = '  +5.42%  '

$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("D36").Value = '''6.15'
$ws.Range("E36").Value = '  +10.73%  '

$ws.Range("D37").Value = '''1.01'
$ws.Range("E37").Value = '  +7.09%  '

$ws.Range("D38").Value = '''47.99'
$ws.Range("E38").Value = '  +1.72%  '

$ws.Range("D39").Value = '''0.330'
$ws.Range("E39").Value = '  +10.31%  '

$ws.Range("D40").Value = '''2.11'
$ws.Range("E40").Value = '  +7.72%  '

$ws.Range("E41").Value = '  +5.18%  '

$ws.Range("D42").Value = '''50.13'
$ws.Range("E42").Value = '  +2.19%  '

$ws.Range("D43").Value = '''8.73'
$ws.Range("E43").Value = '  +5.72%  '

$ws.Range("D44").Value = '''2.83'
$ws.Range("E44").Value = '  +13.10%  '

$ws.Range("D45").Value = '''408.60'
$ws.Range("E45").Value = '  +12.14%  '

$ws.Range("D46").Value = '2.975.03'
$ws.Range("E46").Value = '  +6.83%  '

$ws.Range("D47").Value = '''28.05'
$ws.Range("E47").Value = '  +14.85%  '

$ws.Range("E48").Value = '  +4.71%  '

$ws.Range("D49").Value = '''135.40'
$ws.Range("E49").Value = '  +0.58%  '

$ws.Range("D51").Value = '''2.45'
$ws.Range("E51").Value = '  +13.85%  '
